# Update "想去人数" (number of people interested) counts that changed
# between the two data pulls captured in sheet "展览" and sheet "全部类型".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - rows 3,4,5,7,9,11,21,22,23,27,34,35 in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 336
$ws1.Range("F4").Value  = 424
$ws1.Range("F5").Value  = 1724
$ws1.Range("F7").Value  = 2176
$ws1.Range("F9").Value  = 283
$ws1.Range("F11").Value = 4890
$ws1.Range("F21").Value = 3835
$ws1.Range("F22").Value = 703
$ws1.Range("F23").Value = 647
$ws1.Range("F27").Value = 116
$ws1.Range("F34").Value = 928
$ws1.Range("F35").Value = 2440

# Sheet "全部类型" (fourth sheet) - same events, but shifted down by one row
# starting at row 34 (this sheet has an extra row that "展览" does not have),
# so the last two updates land on F35/F36 instead of F34/F35.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 336
$ws4.Range("F4").Value  = 424
$ws4.Range("F5").Value  = 1724
$ws4.Range("F7").Value  = 2176
$ws4.Range("F9").Value  = 283
$ws4.Range("F11").Value = 4890
$ws4.Range("F21").Value = 3835
$ws4.Range("F22").Value = 703
$ws4.Range("F23").Value = 647
$ws4.Range("F27").Value = 116
$ws4.Range("F35").Value = 928
$ws4.Range("F36").Value = 2440

$wb.Save()
